# Update factsheets with text edits from COMM
#
# The source data previously stored several "count" columns (and a couple of
# "0" placeholder rows) as real numbers. They are converted to literal text
# values instead (so the workbook round-trips every number/percentage/dollar
# figure as a string). A new "Total" row is also appended to the County
# sheet. A leading apostrophe is used where a value would otherwise be
# auto-recognised by Excel as a number/percentage/currency, which forces it
# to be stored as text.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overall": A2 (772) becomes text
# ---------------------------------------------------------------------
$wsOverall = $wb.Worksheets.Item("Overall")
$wsOverall.Range("A2").Value = "'772"

# ---------------------------------------------------------------------
# Sheet "County": B2:B32 counts become text, rows 33/34 get new
# percentage/dollar placeholder text, and a new Total row (35) is added.
# ---------------------------------------------------------------------
$wsCounty = $wb.Worksheets.Item("County")

$countyCounts = @{
    2  = 283
    3  = 10
    4  = 5
    5  = 4
    6  = 9
    7  = 1
    8  = 53
    9  = 20
    10 = 14
    11 = 1
    12 = 1
    13 = 21
    14 = 10
    15 = 10
    16 = 9
    17 = 21
    18 = 2
    19 = 15
    20 = 2
    21 = 17
    22 = 5
    23 = 41
    24 = 8
    25 = 32
    26 = 143
    27 = 5
    28 = 3
    29 = 16
    30 = 1
    31 = 3
    32 = 7
}

foreach ($row in $countyCounts.Keys) {
    $wsCounty.Range("B$row").Value = "'" + $countyCounts[$row]
}

# Catron County (row 33) and Harding County (row 34): replace the "0"
# placeholders with formatted text placeholders.
foreach ($row in 33, 34) {
    $wsCounty.Range("B$row").Value = "'0.00%"
    $wsCounty.Range("C$row").Value = "'`$0"
    $wsCounty.Range("D$row").Value = "'0.00%"
    $wsCounty.Range("E$row").Value = "'0.00%"
    $wsCounty.Range("F$row").Value = "'0.00%"
}

# New Total row at the bottom of the County sheet.
$wsCounty.Range("A35").Value = "Total"
$wsCounty.Range("B35").Value = "'772"
$wsCounty.Range("C35").Value = "'`$943,744,428"
$wsCounty.Range("D35").Value = "'10.15%"
$wsCounty.Range("E35").Value = "'-16.09%"
$wsCounty.Range("F35").Value = "'70.47%"

# ---------------------------------------------------------------------
# Sheet "Congressional District": B2:B5 counts become text
# ---------------------------------------------------------------------
$wsCd = $wb.Worksheets.Item("Congressional District")
$cdCounts = @{ 2 = 301; 3 = 147; 4 = 324; 5 = 772 }
foreach ($row in $cdCounts.Keys) {
    $wsCd.Range("B$row").Value = "'" + $cdCounts[$row]
}

# ---------------------------------------------------------------------
# Sheet "Size": B2:B8 counts become text
# ---------------------------------------------------------------------
$wsSize = $wb.Worksheets.Item("Size")
$sizeCounts = @{ 2 = 275; 3 = 227; 4 = 140; 5 = 38; 6 = 55; 7 = 37; 8 = 772 }
foreach ($row in $sizeCounts.Keys) {
    $wsSize.Range("B$row").Value = "'" + $sizeCounts[$row]
}

# ---------------------------------------------------------------------
# Sheet "Subsector": B2:B13 counts become text
# ---------------------------------------------------------------------
$wsSub = $wb.Worksheets.Item("Subsector")
$subCounts = @{ 2 = 82; 3 = 79; 4 = 54; 5 = 79; 6 = 8; 7 = 231; 8 = 7; 9 = 72; 10 = 9; 11 = 146; 12 = 5; 13 = 772 }
foreach ($row in $subCounts.Keys) {
    $wsSub.Range("B$row").Value = "'" + $subCounts[$row]
}
